# Add a new "Modelo" column (F) with header + value, matching the style
# already used by the existing header row, and refresh the MSE value in B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1, formatted like the rest of the header row (bold,
# bordered, centered) by copying the format from the adjacent header cell.
$ws.Range("F1").Value = "Modelo"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data cell F2 with the model description.
$ws.Range("F2").Value = "Pipeline(steps=[('model', GradientBoostingRegressor(n_estimators=150))])"

# Updated MSE value (tiny floating point refinement).
$ws.Range("B2").Value = 0.03244348355771105
